# Fruta / hortaliza, semanal
# Updates date (D) and volume/price (M, N, O, P, S) values for the
# "Femacal de La Calera - Breva" weekly price rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44193
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 3000

# Row 3
$ws.Range("D3").Value = 44196
$ws.Range("M3").Value = 56

# Row 4
$ws.Range("D4").Value = 44188
$ws.Range("M4").Value = 30

# Row 7
$ws.Range("D7").Value = 44186
$ws.Range("M7").Value = 40

# Row 8
$ws.Range("D8").Value = 44181
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("S8").Value = 4000
